# Team Project Diary - Week 4 log update
# Adds the newest work-log entries to the "Week 4" sheet.
#
# Entry order matters: it determines the order new strings are appended
# to the shared-string table, so the three "new work" entries are typed
# first (B21:C23), and the meeting note that was added afterwards
# (B20:C20) is typed last - matching how the log was actually filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 4")

$ws.Range("B21").Value = "Email Template (2015-10-24)"
$ws.Range("C21").Value = 1

$ws.Range("B22").Value = "StackOverflow/Github/Twitter API Research (2015-10-24)"
$ws.Range("C22").Value = 3

$ws.Range("B23").Value = "Web Scraper Python Script (2015-10-24)"
$ws.Range("C23").Value = 4

$ws.Range("B20").Value = "Thursday Meeting on Survey Questions (2015-10-22)"
$ws.Range("C20").Value = 2.5

$ws.Range("B22").Select() | Out-Null
